# Refresh the scraped cryptocurrency market data (Price / Volume(1h) columns)
# as produced by the "Updated cryptos list ... with GitHub Actions" job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> @(new Price text, new Volume(1h) text). Use $null to leave a
# column untouched for that row.
$updates = @{
  2  = @("29.202.64",     "  +0.10%  ")
  3  = @("1.835.55",      "  -0.36%  ")
  4  = @("0.9992",        $null)
  5  = @("240.82",        "  -0.19%  ")
  6  = @("0.6667",        "  -2.96%  ")
  7  = @($null,           "  -0.01%  ")
  8  = @("0.07372",       "  -1.12%  ")
  9  = @("0.2923",        "  -2.50%  ")
  10 = @("22.67",         "  -2.35%  ")
  11 = @("0.07715",       "  +0.77%  ")
  12 = @("1.825.62",      "  -0.87%  ")
  13 = @("4.979",         "  -1.54%  ")
  14 = @("0.6664",        "  -2.29%  ")
  15 = @("83.45",         "  -4.28%  ")
  16 = @("6.092",         "  -1.11%  ")
  17 = @("29.149.38",     "  -0.09%  ")
  18 = @("0.000008270",   "  +1.10%  ")
  19 = @("225.95",        "  -1.27%  ")
  20 = @("12.44",         "  -0.77%  ")
  21 = @("1.000",         "  +0.07%  ")
  22 = @("7.129",         "  -3.60%  ")
  23 = @("0.9998",        "  -0.03%  ")
  24 = @("160.63",        "  +0.73%  ")
  25 = @("8.637",         "  -1.37%  ")
  26 = @("0.1392",        "  -3.85%  ")
  27 = @("17.92",         "  -0.86%  ")
  28 = @("1.505",         "  -0.57%  ")
  29 = @("4.112",         "  -3.93%  ")
  30 = @("4.032",         "  -2.64%  ")
  31 = @("1.185",         "  -1.00%  ")
  32 = @("0.05301",       "  +0.65%  ")
  33 = @($null,           "  +0.98%  ")
  34 = @("0.7536",        "  -1.02%  ")
  35 = @("1.131",         "  -0.38%  ")
  36 = @($null,           "  -0.32%  ")
  37 = @("1.301.87",      "  -0.32%  ")
  38 = @("0.01797",       "  -1.89%  ")
  39 = @($null,           "  -0.09%  ")
  40 = @("0.9220",        "  -1.66%  ")
  41 = @("0.08618",       "  +15.79%  ")
  42 = @("5.948",         "  -0.12%  ")
  44 = @("102.31",        "  -2.46%  ")
  45 = @($null,           "  +3.88%  ")
  46 = @("1.972.78",      "  -0.75%  ")
  47 = @($null,           "  -0.66%  ")
  48 = @("1.766",         "  -0.18%  ")
  49 = @("63.29",         "  -2.42%  ")
}

foreach ($row in $updates.Keys) {
  $pair = $updates[$row]
  $newPrice = $pair[0]
  $newVolume = $pair[1]

  if ($newPrice -ne $null) {
    # Force the Price cell to remain plain text (matching the workbook's
    # existing convention) instead of letting Excel auto-coerce it to a
    # number, then restore the default "Normal" style so no stray number
    # format is left behind.
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $newPrice
    $cell.Style = "Normal"
  }

  if ($newVolume -ne $null) {
    $ws.Cells.Item($row, 5).Value = $newVolume
  }
}

# Rows 50 and 51 swapped ranking order: EnergySwap now appears above Cronos,
# and both coins received refreshed price/volume figures.
$cellD50 = $ws.Cells.Item(50, 4)
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cellD50.NumberFormat = "@"
$cellD50.Value = "9.048"
$cellD50.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -4.64%  "

$cellD51 = $ws.Cells.Item(51, 4)
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cellD51.NumberFormat = "@"
$cellD51.Value = "0.05930"
$cellD51.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.42%  "
